$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 - copy formatting (bold font, border, centered) from an
# existing header cell, then set its text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F1").Value = "Status"

# Row 2 - replaced summary labels/values
$ws.Range("A2").Value = "Total comparisons:"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "PyType Wins:"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "Scalpel Wins:"
$ws.Range("F2").Value = 0

# Row 3 - replaced summary labels/values (A3/B3 stay blank, as before)
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Accuracy over PyType"
$ws.Range("F3").Value = 100

# New fill style (solid white) applied to the whole summary block A2:F3
$ws.Range("A2:F3").Interior.Color = 16777215
